# New test case added: append a new demand row (row 3) to Sheet1,
# clear the now-blank "hiringplan" value on row 2, and move the
# viewport/selection the way the author left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "hiringplan" (D2) was blanked out, format kept ---
$ws.Range("D2").ClearContents()

# --- Row 3: brand new demand record ---
$ws.Range("A3").Value = "QA"
$ws.Range("B3").Value = "Senior Consultant-Automation"

# "positions" must be stored as text ("3"), matching column C's
# existing Text number format (same style as C2).
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "3"

$ws.Range("D3").Value = "Project Hire"
$ws.Range("E3").Value = "Project Test 4"
$ws.Range("F3").Value = "Regular"
$ws.Range("G3").Value = "India"
$ws.Range("H3").Value = "Chennai"
$ws.Range("I3").Value = "C:\Users\KiranPatil\GITHUB WORKSPACE\ERS_Automation_Team2\TestData\samplepdf.pdf"
$ws.Range("M3").Value = "P1"

# Match row 3's formatting exactly as left by the author: B3/D3/E3 picked
# up the "message" style from N2, while F3/G3/M3 match the plain label
# style used elsewhere in row 2, and C3/H3 match their row-2 counterparts.
$ws.Range("N2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("M3").PasteSpecial(-4122)

# --- View state: scrolled right a bit further and landed on L9 ---
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("G1")
$ws.Range("L9").Select()
